$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormatLocal = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "307.75"
Set-TextValue $ws.Range("E2") "-2.42%"
Set-TextValue $ws.Range("D3") "37.31"
Set-TextValue $ws.Range("E3") "-5.53%"
Set-TextValue $ws.Range("D4") "5.095"
Set-TextValue $ws.Range("E4") "-0.53%"
Set-TextValue $ws.Range("D5") "0.07881"
Set-TextValue $ws.Range("E5") "-3.84%"
Set-TextValue $ws.Range("D6") "1.981"
Set-TextValue $ws.Range("E6") "0.87%"
Set-TextValue $ws.Range("D7") "4.325"
Set-TextValue $ws.Range("E7") "1.88%"
Set-TextValue $ws.Range("D8") "8.198"
Set-TextValue $ws.Range("E8") "-0.38%"
Set-TextValue $ws.Range("D9") "3.094"
Set-TextValue $ws.Range("E9") "-4.19%"
Set-TextValue $ws.Range("D10") "0.9242"
Set-TextValue $ws.Range("E10") "-0.46%"
Set-TextValue $ws.Range("D11") "0.1297"
Set-TextValue $ws.Range("E11") "-8.34%"
Set-TextValue $ws.Range("D12") "0.1866"
Set-TextValue $ws.Range("E12") "-5.70%"
Set-TextValue $ws.Range("D13") "0.08687"
Set-TextValue $ws.Range("E13") "-4.60%"
Set-TextValue $ws.Range("D14") "0.03434"
Set-TextValue $ws.Range("E14") "-2.28%"
Set-TextValue $ws.Range("D15") "0.09739"
Set-TextValue $ws.Range("E15") "-0.81%"
Set-TextValue $ws.Range("D16") "0.001392"
Set-TextValue $ws.Range("E16") "-0.68%"
Set-TextValue $ws.Range("D17") "0.005939"
Set-TextValue $ws.Range("E17") "0.26%"
Set-TextValue $ws.Range("E18") "1,777.70%"
Set-TextValue $ws.Range("D19") "3.573"
Set-TextValue $ws.Range("E19") "-2.12%"
Set-TextValue $ws.Range("D20") "0.3439"
Set-TextValue $ws.Range("E20") "-0.51%"
Set-TextValue $ws.Range("D21") "0.1285"
Set-TextValue $ws.Range("E21") "-1.40%"
Set-TextValue $ws.Range("D22") "5.022"
Set-TextValue $ws.Range("E22") "4.18%"
Set-TextValue $ws.Range("D23") "0.2501"
Set-TextValue $ws.Range("E23") "2.15%"
Set-TextValue $ws.Range("D24") "0.04331"
Set-TextValue $ws.Range("E24") "-1.05%"
Set-TextValue $ws.Range("D25") "0.001224"
Set-TextValue $ws.Range("E25") "0.20%"
Set-TextValue $ws.Range("D26") "0.004598"
Set-TextValue $ws.Range("E26") "-3.73%"
Set-TextValue $ws.Range("E27") "176.96%"
Set-TextValue $ws.Range("D39") "0.02291"
Set-TextValue $ws.Range("E39") "4.02%"
Set-TextValue $ws.Range("D40") "0.04983"
Set-TextValue $ws.Range("E40") "-3.49%"
Set-TextValue $ws.Range("D41") "0.007503"
Set-TextValue $ws.Range("E41") "-0.52%"
Set-TextValue $ws.Range("D42") "0.009978"
Set-TextValue $ws.Range("E42") "1.60%"
Set-TextValue $ws.Range("D43") "0.1352"
Set-TextValue $ws.Range("E43") "-1.57%"
Set-TextValue $ws.Range("D44") "0.002099"
Set-TextValue $ws.Range("E44") "-0.39%"
Set-TextValue $ws.Range("D45") "0.008038"
Set-TextValue $ws.Range("E45") "-18.07%"
Set-TextValue $ws.Range("D46") "0.00006370"
Set-TextValue $ws.Range("E46") "-0.18%"
Set-TextValue $ws.Range("E47") "0.46%"
Set-TextValue $ws.Range("D48") "0.003009"
Set-TextValue $ws.Range("E48") "8.77%"
Set-TextValue $ws.Range("D49") "0.001205"
Set-TextValue $ws.Range("E49") "0.46%"
Set-TextValue $ws.Range("D50") "0.00002109"
Set-TextValue $ws.Range("E50") "0.46%"
Set-TextValue $ws.Range("D51") "0.0002008"
Set-TextValue $ws.Range("E51") "0.46%"
